$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new numeric-looking values would
# otherwise be auto-converted to numbers by Excel, so they remain text
# (matching the source data which is stored as text).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'

# Apply the updated cell values from the latest crypto data pull.
$ws.Range('D2').Value = '67.813.22'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '2.501.31'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '587.75'
$ws.Range('E5').Value = '  +1.18%  '
$ws.Range('D6').Value = '176.13'
$ws.Range('E6').Value = '  +4.29%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '0.517'
$ws.Range('E8').Value = '  +1.32%  '
$ws.Range('E9').Value = '  +5.34%  '
$ws.Range('E10').Value = '  +1.10%  '
$ws.Range('E11').Value = '  +3.94%  '
$ws.Range('E12').Value = '  +1.44%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.943.81'
$ws.Range('E13').Value = '  +1.39%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = '25.79'
$ws.Range('E14').Value = '  +2.68%  '
$ws.Range('D15').Value = '67.716.61'
$ws.Range('E15').Value = '  +1.40%  '
$ws.Range('E16').Value = '  +2.99%  '
$ws.Range('D17').Value = '2.493.66'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').Value = '11.11'
$ws.Range('E18').Value = '  +2.28%  '
$ws.Range('D19').Value = '7.49'
$ws.Range('E19').Value = '  +2.23%  '
$ws.Range('D20').Value = '352.65'
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('D21').Value = '4.09'
$ws.Range('E21').Value = '  +2.42%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').Value = '70.76'
$ws.Range('E23').Value = '  +3.26%  '
$ws.Range('D24').Value = '4.26'
$ws.Range('E24').Value = '  +2.14%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').Value = '9.24'
$ws.Range('E26').Value = '  +2.01%  '
$ws.Range('D27').Value = '2.624.95'
$ws.Range('E27').Value = '  +1.61%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('D29').Value = '0.0₃0916'
$ws.Range('E29').Value = '  +2.95%  '
$ws.Range('D30').Value = '514.77'
$ws.Range('E30').Value = '  +1.03%  '
$ws.Range('D31').Value = '7.86'
$ws.Range('E31').Value = '  +4.07%  '
$ws.Range('D32').Value = '1.26'
$ws.Range('E32').Value = '  +3.78%  '
$ws.Range('E33').Value = '  +1.95%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E35').Value = '  +7.95%  '
$ws.Range('D36').Value = '161.76'
$ws.Range('E36').Value = '  +2.12%  '
$ws.Range('D37').Value = '18.49'
$ws.Range('E37').Value = '  +2.05%  '
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('E39').Value = '  +2.02%  '
$ws.Range('D40').Value = '1.77'
$ws.Range('E40').Value = '  +6.62%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  +2.88%  '
$ws.Range('D43').Value = '4.88'
$ws.Range('E43').Value = '  +3.19%  '
$ws.Range('E44').Value = '  +4.32%  '
$ws.Range('D45').Value = '145.28'
$ws.Range('E45').Value = '  +3.33%  '
$ws.Range('E46').Value = '  +3.54%  '
$ws.Range('E47').Value = '  +4.51%  '
$ws.Range('D48').Value = '0.517'
$ws.Range('E48').Value = '  +2.16%  '
$ws.Range('D49').Value = '0.0747'
$ws.Range('E49').Value = '  +2.75%  '
$ws.Range('E50').Value = '  +2.86%  '
$ws.Range('E51').Value = '  +1.52%  '
